$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12.45259105729622
$ws.Range("C2").Value = 9.842703261832062
$ws.Range("D2").Value = 5.956106969087947
$ws.Range("E2").Value = 16.60502089598162
$ws.Range("G2").Value = 20.93877928243117
$ws.Range("H2").Value = 11.84167763901382
$ws.Range("I2").Value = 16.07648682201348
$ws.Range("N2").Value = 15.51855463649094
$ws.Range("O2").Value = 17.02873520852619
# Row 3
$ws.Range("B3").Value = 11.78067506210959
$ws.Range("C3").Value = 9.316828665898219
$ws.Range("D3").Value = 5.831789150564728
$ws.Range("E3").Value = 15.65593839455784
$ws.Range("G3").Value = 20.66430610313226
$ws.Range("H3").Value = 11.87014552292492
$ws.Range("I3").Value = 16.16736384541136
$ws.Range("N3").Value = 15.55265622031093
$ws.Range("O3").Value = 17.02280610789565
# Row 4
$ws.Range("B4").Value = 11.34836711799554
$ws.Range("C4").Value = 8.976595823729708
$ws.Range("D4").Value = 5.755894916635675
$ws.Range("E4").Value = 15.04774331270542
$ws.Range("G4").Value = 20.50459626066835
$ws.Range("H4").Value = 11.89043252078503
$ws.Range("I4").Value = 16.22772418654024
$ws.Range("N4").Value = 15.5753219973269
$ws.Range("O4").Value = 17.02491259233098
# Row 5
$ws.Range("B5").Value = 11.16739676597386
$ws.Range("C5").Value = 8.83366010069407
$ws.Range("D5").Value = 5.725126245671136
$ws.Range("E5").Value = 14.79376481470514
$ws.Range("G5").Value = 20.44182278097064
$ws.Range("H5").Value = 11.89940339836621
$ws.Range("I5").Value = 16.25346377865875
$ws.Range("N5").Value = 15.5849934853182
$ws.Range("O5").Value = 17.02721308293745
# Row 6
$ws.Range("B6").Value = 11.13706223331008
$ws.Range("C6").Value = 8.809669193710915
$ws.Range("D6").Value = 5.72002816601874
$ws.Range("E6").Value = 14.75123024037612
$ws.Range("G6").Value = 20.43154134691738
$ws.Range("H6").Value = 11.90093544312136
$ws.Range("I6").Value = 16.25780667015595
$ws.Range("N6").Value = 15.58662572115425
$ws.Range("O6").Value = 17.02768206147755
# Row 7
$ws.Range("B7").Value = 11.34594568802912
$ws.Range("C7").Value = 8.974685396651475
$ws.Range("D7").Value = 5.755479251349746
$ws.Range("E7").Value = 15.04434251497534
$ws.Range("G7").Value = 20.50374020740657
$ws.Range("H7").Value = 11.89055065876848
$ws.Range("I7").Value = 16.22806670099821
$ws.Range("N7").Value = 15.57545066822639
$ws.Range("O7").Value = 17.02493778366723
# Row 8
$ws.Range("B8").Value = 12.22510256129667
$ws.Range("C8").Value = 9.665020585944543
$ws.Range("D8").Value = 5.913180466555808
$ws.Range("E8").Value = 16.28320465616874
$ws.Range("G8").Value = 20.84236819320843
$ws.Range("H8").Value = 11.85090934609924
$ws.Range("I8").Value = 16.10687160969751
$ws.Range("N8").Value = 15.529954764759
$ws.Range("O8").Value = 17.0254966283435
# Row 9
$ws.Range("B9").Value = 13.78666251411019
$ws.Range("C9").Value = 10.87866219674897
$ws.Range("D9").Value = 6.223804516723461
$ws.Range("E9").Value = 18.62244081700415
$ws.Range("G9").Value = 21.57201478349177
$ws.Range("H9").Value = 11.79554429394902
$ws.Range("I9").Value = 15.90561641751734
$ws.Range("N9").Value = 15.45441371474248
$ws.Range("O9").Value = 17.07226901834259
# Row 10
$ws.Range("B10").Value = 14.82886789669095
$ws.Range("C10").Value = 11.68240967768278
$ws.Range("D10").Value = 6.450147289272349
$ws.Range("E10").Value = 20.27898379871882
$ws.Range("G10").Value = 22.14188519452335
$ws.Range("H10").Value = 11.76862657240399
$ws.Range("I10").Value = 15.78024639288336
$ws.Range("N10").Value = 15.4072129491862
$ws.Range("O10").Value = 17.13449681639049
# Row 11
$ws.Range("B11").Value = 15.27929592248336
$ws.Range("C11").Value = 12.02862440457294
$ws.Range("D11").Value = 6.552173822419537
$ws.Range("E11").Value = 20.9900733905435
$ws.Range("G11").Value = 22.40709558416966
$ws.Range("H11").Value = 11.75939112611036
$ws.Range("I11").Value = 15.72816127404673
$ws.Range("N11").Value = 15.38753435927191
$ws.Range("O11").Value = 17.16883123613025
# Row 12
$ws.Range("B12").Value = 15.44640247704637
$ws.Range("C12").Value = 12.15691555007192
$ws.Range("D12").Value = 6.590632187705723
$ws.Range("E12").Value = 21.25327283269063
$ws.Range("G12").Value = 22.50826634546753
$ws.Range("H12").Value = 11.75632820601415
$ws.Range("I12").Value = 15.70915479669181
$ws.Range("N12").Value = 15.38033985909519
$ws.Range("O12").Value = 17.18269497478835
# Row 13
$ws.Range("B13").Value = 15.41056777470941
$ws.Range("C13").Value = 12.12941113538922
$ws.Range("D13").Value = 6.582357998604619
$ws.Range("E13").Value = 21.19685776348977
$ws.Range("G13").Value = 22.48644630706659
$ws.Range("H13").Value = 11.75696852017679
$ws.Range("I13").Value = 15.7132161995615
$ws.Range("N13").Value = 15.38187788561482
$ws.Range("O13").Value = 17.17967091928958
# Row 14
$ws.Range("B14").Value = 15.2931135166236
$ws.Range("C14").Value = 12.03923546900884
$ws.Range("D14").Value = 6.555341591525046
$ws.Range("E14").Value = 21.01184852115797
$ws.Range("G14").Value = 22.41540477778569
$ws.Range("H14").Value = 11.75913042464144
$ws.Range("I14").Value = 15.72658319630764
$ws.Range("N14").Value = 15.38693730855966
$ws.Range("O14").Value = 17.1699545690232
# Row 15
$ws.Range("B15").Value = 15.22071724080378
$ws.Range("C15").Value = 11.98363342303385
$ws.Range("D15").Value = 6.538769025394892
$ws.Range("E15").Value = 20.89773490629673
$ws.Range("G15").Value = 22.37198285203444
$ws.Range("H15").Value = 11.76051126115182
$ws.Range("I15").Value = 15.73486441739485
$ws.Range("N15").Value = 15.39006985124657
$ws.Range("O15").Value = 17.16411512649496
# Row 16
$ws.Range("B16").Value = 14.79894985407603
$ws.Range("C16").Value = 11.65939112956741
$ws.Range("D16").Value = 6.443457094623099
$ws.Range("E16").Value = 20.23166051054007
$ws.Range("G16").Value = 22.12466320931076
$ws.Range("H16").Value = 11.76929082937428
$ws.Range("I16").Value = 15.78375030104624
$ws.Range("N16").Value = 15.40853503224276
$ws.Range("O16").Value = 17.13237389349339
# Row 17
$ws.Range("B17").Value = 14.53410116715687
$ws.Range("C17").Value = 11.45548964140566
$ws.Range("D17").Value = 6.384715413800089
$ws.Range("E17").Value = 19.81218941232551
$ws.Range("G17").Value = 21.97438679810698
$ws.Range("H17").Value = 11.77544881037528
$ws.Range("I17").Value = 15.81501115572794
$ws.Range("N17").Value = 15.42032175296879
$ws.Range("O17").Value = 17.11444243731378
# Row 18
$ws.Range("B18").Value = 14.37954115076339
$ws.Range("C18").Value = 11.33638440833631
$ws.Range("D18").Value = 6.350842827387308
$ws.Range("E18").Value = 19.56692029329695
$ws.Range("G18").Value = 21.88852109025602
$ws.Range("H18").Value = 11.77927391016399
$ws.Range("I18").Value = 15.83345679728518
$ws.Range("N18").Value = 15.42726998462734
$ws.Range("O18").Value = 17.10469636299942
# Row 19
$ws.Range("B19").Value = 14.3268293732094
$ws.Range("C19").Value = 11.2957444628042
$ws.Range("D19").Value = 6.339360714509079
$ws.Range("E19").Value = 19.48318826998582
$ws.Range("G19").Value = 21.85954980567455
$ws.Range("H19").Value = 11.78061761070757
$ws.Range("I19").Value = 15.83978190314152
$ws.Range("N19").Value = 15.42965154767777
$ws.Range("O19").Value = 17.10149411655771
# Row 20
$ws.Range("B20").Value = 14.56252559602434
$ws.Range("C20").Value = 11.47738450933224
$ws.Range("D20").Value = 6.390977771156749
$ws.Range("E20").Value = 19.85725661140079
$ws.Range("G20").Value = 21.9903259106516
$ws.Range("H20").Value = 11.77476396144029
$ws.Range("I20").Value = 15.81163519383862
$ws.Range("N20").Value = 15.4190495677731
$ws.Range("O20").Value = 17.11629255307798
# Row 21
$ws.Range("B21").Value = 15.32770699076428
$ws.Range("C21").Value = 12.06579872846923
$ws.Range("D21").Value = 6.563282087456327
$ws.Range("E21").Value = 21.06635481885261
$ws.Range("G21").Value = 22.43625221188587
$ws.Range("H21").Value = 11.7584836213883
$ws.Range("I21").Value = 15.72263747934831
$ws.Range("N21").Value = 15.38544425292459
$ws.Range("O21").Value = 17.17278514471373
# Row 22
$ws.Range("B22").Value = 15.80760737593471
$ws.Range("C22").Value = 12.43395996953183
$ws.Range("D22").Value = 6.674844322474377
$ws.Range("E22").Value = 21.82117762151871
$ws.Range("G22").Value = 22.73195799988705
$ws.Range("H22").Value = 11.75037575160867
$ws.Range("I22").Value = 15.66865441378769
$ws.Range("N22").Value = 15.36498106527578
$ws.Range("O22").Value = 17.21472782111734
# Row 23
$ws.Range("B23").Value = 15.55333838948734
$ws.Range("C23").Value = 12.23897167457737
$ws.Range("D23").Value = 6.615410619424999
$ws.Range("E23").Value = 21.42154180590184
$ws.Range("G23").Value = 22.57378247642988
$ws.Range("H23").Value = 11.75447090250521
$ws.Range("I23").Value = 15.69708161736706
$ws.Range("N23").Value = 15.37576558708895
$ws.Range("O23").Value = 17.19188463798404
# Row 24
$ws.Range("B24").Value = 14.54968205569185
$ws.Range("C24").Value = 11.4674916923801
$ws.Range("D24").Value = 6.3881468714479
$ws.Range("E24").Value = 19.83689453445091
$ws.Range("G24").Value = 21.98311817726457
$ws.Range("H24").Value = 11.77507269458999
$ws.Range("I24").Value = 15.81315999207655
$ws.Range("N24").Value = 15.41962418720774
$ws.Range("O24").Value = 17.11545436174204
# Row 25
$ws.Range("B25").Value = 13.38231603209362
$ws.Range("C25").Value = 10.56563656806774
$ws.Range("D25").Value = 6.139912168608453
$ws.Range("E25").Value = 17.97450847294245
$ws.Range("G25").Value = 21.36822862306475
$ws.Range("H25").Value = 11.80811386502455
$ws.Range("I25").Value = 16.07648682201348
$ws.Range("N25").Value = 15.47338942186599
$ws.Range("O25").Value = 17.05471845060718
